# "Cleaned up Fish XLS file and metadata"
#
# - Rename "Sheet1" -> "Master"
# - Remove the obsolete "Sheet2" worksheet (its only unique content, the
#   "Common" header string, was a duplicate of information already present
#   on the Piedmont/MACP sheets)
# - Make "Master" the active/selected tab (was previously "Piedmont")
#
# Note: the workbook-level x15ac:absPath (last-saved-from folder hint) isn't
# reachable through the Excel object model (Path/FullName are read-only
# no-ops for that MS extension attribute in this runtime), so it can't be
# updated from COM automation.

$wb = $excel.ActiveWorkbook

# Rename the primary data sheet.
$master = $wb.Worksheets("Sheet1")
$master.Name = "Master"

# Drop the redundant "Sheet2" worksheet entirely (superseded by Piedmont).
$wb.Worksheets("Sheet2").Delete() | Out-Null

# Select "Master" so it becomes the workbook's active tab.
$master.Select() | Out-Null

Write-Output "Renamed Sheet1->Master, removed Sheet2, selected Master."
